$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every data value (even numeric-looking GasUsed figures) as
# text, matching the workbook author's original export format. Pre-format the
# rows we touch (existing rows 2-58 plus the two newly appended rows 59-60) as
# Text so COM assignment does not silently coerce numeric strings to numbers.
$ws.Range("A2:F60").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(7, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(7, 5).Value = '91010'
$ws.Cells.Item(7, 6).Value = '0.0009101'
$ws.Cells.Item(8, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(9, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(9, 5).Value = '100363'
$ws.Cells.Item(9, 6).Value = '0.00100363'
$ws.Cells.Item(10, 5).Value = '4862386'
$ws.Cells.Item(10, 6).Value = '0.04862386'
$ws.Cells.Item(15, 5).Value = '429683'
$ws.Cells.Item(15, 6).Value = '0.00429683'
$ws.Cells.Item(16, 5).Value = '437530'
$ws.Cells.Item(16, 6).Value = '0.0043753'
$ws.Cells.Item(17, 5).Value = '482606'
$ws.Cells.Item(17, 6).Value = '0.00482606'
$ws.Cells.Item(18, 5).Value = '588145'
$ws.Cells.Item(18, 6).Value = '0.00588145'
$ws.Cells.Item(19, 5).Value = '340683'
$ws.Cells.Item(19, 6).Value = '0.00340683'
$ws.Cells.Item(30, 5).Value = '480828'
$ws.Cells.Item(30, 6).Value = '0.00480828'
$ws.Cells.Item(31, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(32, 1).Value = 'DepositManager'
$ws.Cells.Item(32, 2).Value = 'approveAndCall'
$ws.Cells.Item(32, 4).Value = 'deposit TON to DAO Candidate'
$ws.Cells.Item(32, 5).Value = '340614'
$ws.Cells.Item(32, 6).Value = '0.00340614'
$ws.Cells.Item(33, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(33, 2).Value = 'registerLayer2Candidate'
$ws.Cells.Item(33, 4).Value = 'thanos SystemConfig'
$ws.Cells.Item(33, 5).Value = '4823495'
$ws.Cells.Item(33, 6).Value = '0.04823495'
$ws.Cells.Item(34, 1).Value = 'SeigManager'
$ws.Cells.Item(34, 2).Value = 'updateSeigniorageLayer'
$ws.Cells.Item(34, 4).Value = 'titanLayerAddress'
$ws.Cells.Item(34, 5).Value = '355945'
$ws.Cells.Item(34, 6).Value = '0.00355945'
$ws.Cells.Item(35, 4).Value = 'with claim'
$ws.Cells.Item(35, 5).Value = '355880'
$ws.Cells.Item(35, 6).Value = '0.0035588'
$ws.Cells.Item(36, 1).Value = 'titanLayerContract'
$ws.Cells.Item(36, 2).Value = 'updateSeigniorage'
$ws.Cells.Item(36, 4).Value = 'with staking'
$ws.Cells.Item(36, 5).Value = '355880'
$ws.Cells.Item(36, 6).Value = '0.0035588'
$ws.Cells.Item(37, 1).Value = 'DepositManager'
$ws.Cells.Item(37, 2).Value = 'deposit(address,uint256)'
$ws.Cells.Item(37, 4).Value = ''
$ws.Cells.Item(37, 5).Value = '286891'
$ws.Cells.Item(37, 6).Value = '0.00286891'
$ws.Cells.Item(38, 1).Value = 'SeigManager'
$ws.Cells.Item(38, 2).Value = 'updateSeigniorageLayer'
$ws.Cells.Item(38, 4).Value = 'no give seigniorage to l2'
$ws.Cells.Item(38, 5).Value = '435149'
$ws.Cells.Item(38, 6).Value = '0.00435149'
$ws.Cells.Item(39, 1).Value = 'DepositManager'
$ws.Cells.Item(39, 2).Value = 'deposit(address,address,uint256)'
$ws.Cells.Item(39, 4).Value = ''
$ws.Cells.Item(39, 5).Value = '254169'
$ws.Cells.Item(39, 6).Value = '0.00254169'
$ws.Cells.Item(40, 4).Value = 'give seigniorage to l2'
$ws.Cells.Item(40, 5).Value = '408534'
$ws.Cells.Item(40, 6).Value = '0.00408534'
$ws.Cells.Item(41, 1).Value = 'SeigManager'
$ws.Cells.Item(41, 2).Value = 'updateSeigniorageLayer'
$ws.Cells.Item(41, 4).Value = 'not operator'
$ws.Cells.Item(41, 5).Value = '408534'
$ws.Cells.Item(41, 6).Value = '0.00408534'
$ws.Cells.Item(42, 4).Value = 'operator'
$ws.Cells.Item(42, 5).Value = '487810'
$ws.Cells.Item(42, 6).Value = '0.0048781'
$ws.Cells.Item(43, 1).Value = 'Layer2Contract'
$ws.Cells.Item(43, 2).Value = 'updateSeigniorage'
$ws.Cells.Item(43, 4).Value = 'with operator''s staking'
$ws.Cells.Item(43, 5).Value = '593349'
$ws.Cells.Item(43, 6).Value = '0.00593349'
$ws.Cells.Item(44, 2).Value = 'requestWithdrawal'
$ws.Cells.Item(44, 5).Value = '341275'
$ws.Cells.Item(44, 6).Value = '0.00341275'
$ws.Cells.Item(45, 1).Value = 'DepositManager'
$ws.Cells.Item(45, 2).Value = 'processRequest'
$ws.Cells.Item(45, 5).Value = '180262'
$ws.Cells.Item(45, 6).Value = '0.00180262'
$ws.Cells.Item(46, 1).Value = 'L1BridgeRegistry'
$ws.Cells.Item(46, 2).Value = 'restoreLayer2Candidate'
$ws.Cells.Item(46, 5).Value = '72294'
$ws.Cells.Item(46, 6).Value = '0.00072294'
$ws.Cells.Item(47, 5).Value = '426030'
$ws.Cells.Item(47, 6).Value = '0.0042603'
$ws.Cells.Item(48, 5).Value = '403330'
$ws.Cells.Item(48, 6).Value = '0.0040333'
$ws.Cells.Item(49, 1).Value = 'SeigManager'
$ws.Cells.Item(49, 2).Value = 'updateSeigniorageLayer'
$ws.Cells.Item(49, 4).Value = ''
$ws.Cells.Item(49, 5).Value = '403330'
$ws.Cells.Item(49, 6).Value = '0.0040333'
$ws.Cells.Item(50, 4).Value = 'with operator''s claim'
$ws.Cells.Item(50, 5).Value = '465506'
$ws.Cells.Item(50, 6).Value = '0.00465506'
$ws.Cells.Item(51, 1).Value = 'Layer2Contract'
$ws.Cells.Item(51, 2).Value = 'updateSeigniorage'
$ws.Cells.Item(51, 4).Value = 'with operator''s staking'
$ws.Cells.Item(51, 5).Value = '567828'
$ws.Cells.Item(51, 6).Value = '0.00567828'
$ws.Cells.Item(52, 1).Value = 'Layer2Contract'
$ws.Cells.Item(52, 2).Value = 'updateSeigniorage'
$ws.Cells.Item(52, 4).Value = 'with operator''s staking'
$ws.Cells.Item(52, 5).Value = '437710'
$ws.Cells.Item(52, 6).Value = '0.0043771'
$ws.Cells.Item(53, 1).Value = 'DepositManager'
$ws.Cells.Item(53, 2).Value = 'requestWithdrawal'
$ws.Cells.Item(53, 4).Value = ''
$ws.Cells.Item(53, 5).Value = '324175'
$ws.Cells.Item(53, 6).Value = '0.00324175'
$ws.Cells.Item(54, 2).Value = 'processRequest'
$ws.Cells.Item(54, 5).Value = '128962'
$ws.Cells.Item(54, 6).Value = '0.00128962'
$ws.Cells.Item(55, 1).Value = 'TonContract'
$ws.Cells.Item(55, 2).Value = 'approveAndCall'
$ws.Cells.Item(55, 4).Value = 'DepositManager.onApprove'
$ws.Cells.Item(55, 5).Value = '322929'
$ws.Cells.Item(55, 6).Value = '0.00322929'
$ws.Cells.Item(56, 1).Value = 'DepositManager'
$ws.Cells.Item(56, 2).Value = 'deposit(address,uint256)'
$ws.Cells.Item(56, 5).Value = '249099'
$ws.Cells.Item(56, 6).Value = '0.00249099'
$ws.Cells.Item(57, 2).Value = 'deposit(address,address,uint256)'
$ws.Cells.Item(57, 5).Value = '257425'
$ws.Cells.Item(57, 6).Value = '0.00257425'
$ws.Cells.Item(58, 1).Value = 'SeigManager'
$ws.Cells.Item(58, 2).Value = 'updateSeigniorageLayer'
$ws.Cells.Item(58, 5).Value = '354699'
$ws.Cells.Item(58, 6).Value = '0.00354699'
$ws.Cells.Item(59, 1).Value = 'DepositManager'
$ws.Cells.Item(59, 2).Value = 'requestWithdrawal'
$ws.Cells.Item(59, 3).Value = ''
$ws.Cells.Item(59, 4).Value = ''
$ws.Cells.Item(59, 5).Value = '326832'
$ws.Cells.Item(59, 6).Value = '0.00326832'
$ws.Cells.Item(60, 1).Value = 'DepositManager'
$ws.Cells.Item(60, 2).Value = 'processRequests'
$ws.Cells.Item(60, 3).Value = ''
$ws.Cells.Item(60, 4).Value = ''
$ws.Cells.Item(60, 5).Value = '109935'
$ws.Cells.Item(60, 6).Value = '0.00109935'

# Extend the ignoredErrors "numbers stored as text" rule and the sheet dimension
# to cover the two newly appended rows (Excel keeps these in sync automatically
# based on the used range, so no further action is required here).
